$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Move the password value from B2 to B4 (Failed UserName/Password screenshot row)
$ws.Range("B2").ClearContents()
$ws.Range("B4").Value = "leo_1"
